# Automatic update of files.
# Re-order the species "records" (taxon info in A/B/D/E/F/G/H together with the
# Ost/Nord coordinates in Q/R) across rows 10-21 and 24-27, and move the
# "overblommad" remark (J/K/L/N/AF) from row 21 to row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---
$ws.Range("A10").Value = 112044189
$ws.Range("Q10").Value = 554686
$ws.Range("R10").Value = 6698721

# --- Row 11 ---
$ws.Range("A11").Value = 112044170
$ws.Range("B11").Value = 89993
$ws.Range("E11").Value = 1209
$ws.Range("F11").Value = "Rynkskinn"
$ws.Range("G11").Value = "Phlebia centrifuga"
$ws.Range("H11").Value = "P.Karst."
$ws.Range("Q11").Value = 554745
$ws.Range("R11").Value = 6698641

# --- Row 12 ---
$ws.Range("A12").Value = 112044187
$ws.Range("Q12").Value = 554629
$ws.Range("R12").Value = 6698775

# --- Row 13 ---
$ws.Range("A13").Value = 112044194
$ws.Range("Q13").Value = 554746
$ws.Range("R13").Value = 6698619

# --- Row 14 ---
$ws.Range("A14").Value = 112044193
$ws.Range("Q14").Value = 554737
$ws.Range("R14").Value = 6698616

# --- Row 15 ---
$ws.Range("A15").Value = 112044188
$ws.Range("Q15").Value = 554647
$ws.Range("R15").Value = 6698760

# --- Row 16 ---
$ws.Range("A16").Value = 112044192
$ws.Range("Q16").Value = 554727
$ws.Range("R16").Value = 6698622

# --- Row 17 --- (gains the "overblommad" remark that used to live on row 21)
$ws.Range("A17").Value = 112044174
$ws.Range("J17").Style = "Normal"
$ws.Range("K17").Value = "överblommad"
$ws.Range("L17").Style = "Normal"
$ws.Range("N17").Style = "Normal"
$ws.Range("Q17").Value = 554690
$ws.Range("R17").Value = 6698722
$ws.Range("AF17").Style = "Normal"

# --- Row 18 ---
$ws.Range("A18").Value = 112044190
$ws.Range("B18").Value = 96735
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("Q18").Value = 554682
$ws.Range("R18").Value = 6698694

# --- Row 19 ---
$ws.Range("A19").Value = 112044191
$ws.Range("Q19").Value = 554719
$ws.Range("R19").Value = 6698669

# --- Row 20 ---
$ws.Range("A20").Value = 112044185
$ws.Range("Q20").Value = 554752
$ws.Range("R20").Value = 6698637

# --- Row 21 --- (loses the "overblommad" remark, now on row 17)
$ws.Range("A21").Value = 112044186
$ws.Range("J21").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("Q21").Value = 554675
$ws.Range("R21").Value = 6698785
$ws.Range("AF21").ClearContents()

# --- Row 24 ---
$ws.Range("A24").Value = 112044158
$ws.Range("B24").Value = 89553
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 1202
$ws.Range("F24").Value = "Ullticka"
$ws.Range("G24").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H24").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q24").Value = 554756
$ws.Range("R24").Value = 6698631

# --- Row 25 ---
$ws.Range("A25").Value = 112044184
$ws.Range("Q25").Value = 554833
$ws.Range("R25").Value = 6698646

# --- Row 26 ---
$ws.Range("A26").Value = 112044171
$ws.Range("B26").Value = 89834
$ws.Range("E26").Value = 658
$ws.Range("F26").Value = "Rosenticka"
$ws.Range("G26").Value = "Rhodofomes roseus"
$ws.Range("H26").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q26").Value = 554758
$ws.Range("R26").Value = 6698625

# --- Row 27 ---
$ws.Range("A27").Value = 112044195
$ws.Range("B27").Value = 96735
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = "Knärot"
$ws.Range("G27").Value = "Goodyera repens"
$ws.Range("H27").Value = "(L.) R. Br."
$ws.Range("Q27").Value = 554806
$ws.Range("R27").Value = 6698598
